$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

# --- Update existing Tech.Debt column (C2:C11) from numeric to descriptive text ---
$ws.Range("C2").Value = '15 min'
$ws.Range("C3").Value = '15 min'
$ws.Range("C4").Value = '15 min'
$ws.Range("C5").Value = '15 min'
$ws.Range("C6").Value = '30 min'
$ws.Range("C7").Value = '15 min'
$ws.Range("C8").Value = '4h'
$ws.Range("C9").Value = '15 min'
$ws.Range("C10").Value = '15 min'
$ws.Range("C11").Value = '15 min'

# Row 2 header-row-adjacent data row grew a touch taller in the real edit
$ws.Rows.Item(2).RowHeight = 30

# --- Append new rule rows 12-26 ---
# Row 12: MQ Node
$ws.Range("A12").Value = 'MQ Node'
$ws.Range("B12").Value = 'MAJOR'
$ws.Range("C12").Value = '15 min'
$ws.Range("D12").Value = 'Bad Practice'
$ws.Range("E12").Value = 'MQ INPUT, MQ OUTPUT, MQ GET : Ensure the transaction mode is set to Automatic to use transactions with persistent messages.'
$ws.Range("E12").WrapText = $true

# Row 13: MQ Node
$ws.Range("A13").Value = 'MQ Node'
$ws.Range("B13").Value = 'MAJOR'
$ws.Range("C13").Value = '15 min'
$ws.Range("D13").Value = 'Bad Practice'
$ws.Range("E13").Value = 'MQ REPLY : Ensure the transaction mode is set to Automatic to use transactions with persistent messages.'
$ws.Range("E13").WrapText = $true

# Row 14: Transformation Node (ComputeNode)
$ws.Range("A14").Value = 'Transformation Node (ComputeNode)'
$ws.Range("B14").Value = 'MAJOR'
$ws.Range("C14").Value = '30 min'
$ws.Range("D14").Value = 'Standards'
$ws.Range("E14").Value = 'Avoid overusing this node as tree copying is processor heavy.  <-- We don''t do this.
Give them verb-noun names in upper-camel case without spaces so  the underlying module can be named normally.  <-- Only do naming convention'
$ws.Range("E14").WrapText = $true

# Row 15: HttpReply Node
$ws.Range("A15").Value = 'HttpReply Node'
$ws.Range("B15").Value = 'MAJOR'
$ws.Range("C15").Value = '30 min'
$ws.Range("D15").Value = 'Bad Practice'
$ws.Range("E15").Value = '"Ignore transport failures" property should be set for HTTPReply Node
'
$ws.Range("E15").WrapText = $true

# Row 16: HttpReply Node
$ws.Range("A16").Value = 'HttpReply Node'
$ws.Range("B16").Value = 'MAJOR'
$ws.Range("C16").Value = '30 min'
$ws.Range("D16").Value = 'Bad Practice'
$ws.Range("E16").Value = '"Generate default HTTP headers from reply or response" property should be set for HTTPReply Node
'
$ws.Range("E16").WrapText = $true

# Row 17: Webservice Node
$ws.Range("A17").Value = 'Webservice Node'
$ws.Range("B17").Value = 'CRITICAL'
$ws.Range("C17").Value = '15 min'
$ws.Range("D17").Value = 'Performance '
$ws.Range("E17").Value = 'Ensure that the "Request timeout" property in the SOAP node is set'
$ws.Range("E17").WrapText = $true

# Row 18: All Nodes
$ws.Range("A18").Value = 'All Nodes'
$ws.Range("B18").Value = 'MAJOR'
$ws.Range("C18").Value = '20 min'
$ws.Range("D18").Value = 'Performance '
$ws.Range("E18").Value = 'Use XMLNSC over XMLNS'
$ws.Range("E18").WrapText = $true

# Row 19: All Nodes
$ws.Range("A19").Value = 'All Nodes'
$ws.Range("B19").Value = 'MAJOR'
$ws.Range("C19").Value = '20 min'
$ws.Range("D19").Value = 'Correctness, readability'
$ws.Range("E19").Value = 'There is no input connection to this node. The code may not be reachable or functioning'
$ws.Range("E19").WrapText = $true

# Row 20: Label Node
$ws.Range("A20").Value = 'Label Node'
$ws.Range("B20").Value = 'MAJOR'
$ws.Range("C20").Value = '30 min'
$ws.Range("D20").Value = 'Correctness'
$ws.Range("E20").Value = 'Label has no associated processing logic attached'
$ws.Range("E20").WrapText = $true

# Row 21: All Nodes
$ws.Range("A21").Value = 'All Nodes'
$ws.Range("B21").Value = 'MAJOR'
$ws.Range("C21").Value = '30 min'
$ws.Range("D21").Value = 'Completeness'
$ws.Range("E21").Value = 'All input terminals are not connected. Processing may not complete normally'
$ws.Range("E21").WrapText = $true

# Row 22: Filter Node
$ws.Range("A22").Value = 'Filter Node'
$ws.Range("B22").Value = 'CRITICAL'
$ws.Range("C22").Value = '1 h'
$ws.Range("D22").Value = 'Correctness'
$ws.Range("E22").Value = 'The filter node may not have its connections connected correctly'
$ws.Range("E22").WrapText = $true

# Row 23: RouteTo Node
$ws.Range("A23").Value = 'RouteTo Node'
$ws.Range("B23").Value = 'MAJOR'
$ws.Range("C23").Value = '20 min'
$ws.Range("D23").Value = 'Correctness, readability'
$ws.Range("E23").Value = 'Usually the RouteTo and Label are in the same flow as to make things more readable'
$ws.Range("E23").WrapText = $true

# Row 24: MessageFlow
$ws.Range("A24").Value = 'MessageFlow'
$ws.Range("B24").Value = 'CRITICAL'
$ws.Range("C24").Value = '30 min'
$ws.Range("D24").Value = 'Correctness'
$ws.Range("E24").Value = 'The message flow does not consistently reply to messages/requests'
$ws.Range("E24").WrapText = $true

# Row 25: Trace Node
$ws.Range("A25").Value = 'Trace Node'
$ws.Range("B25").Value = 'MINOR'
$ws.Range("C25").Value = '5 min'
$ws.Range("D25").Value = 'Standards'
$ws.Range("E25").Value = 'Trace nodes should not be used'
$ws.Range("E25").WrapText = $true

# Row 26: Soap Async Node
$ws.Range("A26").Value = 'Soap Async Node'
$ws.Range("B26").Value = 'CRITICAL'
$ws.Range("C26").Value = '30 min'
$ws.Range("D26").Value = 'Correctness'
$ws.Range("E26").Value = 'The SOAP Async node ''fault'' terminal should be connected'
$ws.Range("E26").WrapText = $true

# Approximate autofit heights for rows whose description wraps onto extra lines
$ws.Rows.Item(14).RowHeight = 45
$ws.Rows.Item(15).RowHeight = 30
$ws.Rows.Item(16).RowHeight = 30
$ws.Rows.Item(19).RowHeight = 30
$ws.Rows.Item(23).RowHeight = 30

# --- Column width adjustments ---
$ws.Columns.Item(4).ColumnWidth = 14.42578125
$ws.Columns.Item(5).ColumnWidth = 128

# --- View / selection ---
$ws.Application.ActiveWindow.ScrollRow = 16
$ws.Range("E21").Select()

# --- Page setup ---
$ws.PageSetup.Orientation = 1